$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" / last-changed date) for rows 2 through 15
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
